$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 501.14285
$ws.Range("I19").Value = 628.25
$ws.Range("J19").Value = 331.66666
$ws.Range("K19").Value = 628.25
$ws.Range("L19").Value = 331.66666
$ws.Range("M19").Value = -453.25
$ws.Range("N19").Value = -681.66666

$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("M21").Value = 10000
$ws.Range("N21").Value = -10936

$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 10000
$ws.Range("N23").Value = -10468

$ws.Range("H34").Value = 4871
$ws.Range("J34").Value = 5999.5
$ws.Range("L34").Value = 5999.5
$ws.Range("N34").Value = -6405.5

$ws.Range("H36").Value = 4871
$ws.Range("J36").Value = 5999.5
$ws.Range("L36").Value = 5999.5
$ws.Range("N36").Value = -7429.5

$ws.Range("H116").Value = 18367
$ws.Range("J116").Value = 8999.75
$ws.Range("L116").Value = 8999.75
$ws.Range("N116").Value = -15883.75

$ws.Range("H137").Value = 85419.39999999999
$ws.Range("I137").Value = 104124.5
$ws.Range("J137").Value = 10599
$ws.Range("K137").Value = 312373.5
$ws.Range("L137").Value = 31797
$ws.Range("M137").Value = -309823.5
$ws.Range("N137").Value = -36897

$ws.Range("H140").Value = 569998.5
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8335983.5
$ws.Range("I32").Value = 3877748.5
$ws.Range("J32").Value = 35722284
$ws.Range("K32").Value = 3877748.5
$ws.Range("L32").Value = 35722284
$ws.Range("M32").Value = -3877461.5
$ws.Range("N32").Value = -35722858

$ws.Range("H80").Value = 100000
$ws.Range("J80").Value = 100000
$ws.Range("L80").Value = 100000
$ws.Range("N80").Value = -101996

$ws.Range("H83").Value = 100000
$ws.Range("J83").Value = 100000
$ws.Range("L83").Value = 300000
$ws.Range("N83").Value = -309984

$ws.Range("H97").Value = 1340.5416
$ws.Range("J97").Value = 2248.5
$ws.Range("L97").Value = 2248.5
$ws.Range("N97").Value = -3240.5

$ws.Range("H132").Value = 2492.06
$ws.Range("I132").Value = 2021.8948
$ws.Range("K132").Value = 6065.6844
$ws.Range("M132").Value = -3535.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11907712
$ws.Range("I134").Value = 3247914.8
$ws.Range("J134").Value = 27784008
$ws.Range("K134").Value = 9743744.399999999
$ws.Range("L134").Value = 83352024
$ws.Range("M134").Value = -9741209.399999999
$ws.Range("N134").Value = -83357094

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2622.7144
$ws.Range("J58").Value = 3624.7778
$ws.Range("L58").Value = 3624.7778
$ws.Range("N58").Value = -4030.7778

$ws.Range("H75").Value = 91499.5
$ws.Range("J75").Value = 102799.4
$ws.Range("L75").Value = 102799.4
$ws.Range("N75").Value = -104795.4

$ws.Range("H78").Value = 91499.5
$ws.Range("J78").Value = 102799.4
$ws.Range("L78").Value = 308398.2
$ws.Range("N78").Value = -318382.2

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0

$ws.Range("H100").Value = 111990
$ws.Range("J100").Value = 111990
$ws.Range("L100").Value = 111990
$ws.Range("N100").Value = -114154

$ws.Range("H132").Value = 1672.375
$ws.Range("J132").Value = 2458.8333
$ws.Range("L132").Value = 7376.499899999999
$ws.Range("N132").Value = -12436.4999

$ws.Range("H136").Value = 2622.7144
$ws.Range("J136").Value = 3624.7778
$ws.Range("L136").Value = 10874.3334
$ws.Range("N136").Value = -15974.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 1823.75
$ws.Range("I124").Value = 1823.75
$ws.Range("K124").Value = 5471.25
$ws.Range("M124").Value = -561.25

$ws.Range("H131").Value = 2311.4546
$ws.Range("I131").Value = 1185.2
$ws.Range("J131").Value = 3250
$ws.Range("K131").Value = 3555.6
$ws.Range("L131").Value = 9750
$ws.Range("M131").Value = 1484.4
$ws.Range("N131").Value = -19830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 84666
$ws.Range("J15").Value = 84666
$ws.Range("L15").Value = 84666
$ws.Range("N15").Value = -85242

$ws.Range("H43").Value = 5633.3335
$ws.Range("I43").Value = 4166.6665
$ws.Range("J43").Value = 7100
$ws.Range("K43").Value = 4166.6665
$ws.Range("L43").Value = 7100
$ws.Range("M43").Value = -4015.6665
$ws.Range("N43").Value = -7402

$ws.Range("H57").Value = 30997.5
$ws.Range("I57").Value = 12000
$ws.Range("K57").Value = 12000
$ws.Range("M57").Value = -11180

$ws.Range("H81").Value = 84666
$ws.Range("J81").Value = 84666
$ws.Range("L81").Value = 84666
$ws.Range("N81").Value = -86662

$ws.Range("H84").Value = 84666
$ws.Range("J84").Value = 84666
$ws.Range("L84").Value = 253998
$ws.Range("N84").Value = -263982

$ws.Range("H102").Value = 2001.5238
$ws.Range("I102").Value = 1696.2222
$ws.Range("J102").Value = 3833.3333
$ws.Range("K102").Value = 1696.2222
$ws.Range("L102").Value = 3833.3333
$ws.Range("M102").Value = -74.22219999999993
$ws.Range("N102").Value = -7077.3333

$ws.Range("H113").Value = 2538
$ws.Range("I113").Value = 2363.8333
$ws.Range("J113").Value = 2799.25
$ws.Range("K113").Value = 2363.8333
$ws.Range("L113").Value = 2799.25
$ws.Range("M113").Value = -193.8332999999998
$ws.Range("N113").Value = -7139.25

$ws.Range("H122").Value = 2993.818
$ws.Range("I122").Value = 2993.818
$ws.Range("K122").Value = 8981.454000000002
$ws.Range("M122").Value = -6531.454000000002

$ws.Range("H126").Value = 2497.8333
$ws.Range("I126").Value = 2197.4
$ws.Range("K126").Value = 6592.200000000001
$ws.Range("M126").Value = -4122.200000000001

$ws.Range("H128").Value = 134989
$ws.Range("J128").Value = 134989
$ws.Range("L128").Value = 134989
$ws.Range("N128").Value = -144949

$ws.Range("H132").Value = 1652.6471
$ws.Range("I132").Value = 1473
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4419
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1889
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 398.07144
$ws.Range("I55").Value = 306.1905
$ws.Range("J55").Value = 673.7143
$ws.Range("K55").Value = 306.1905
$ws.Range("L55").Value = 673.7143
$ws.Range("M55").Value = -133.1905
$ws.Range("N55").Value = -1019.7143

$ws.Range("H122").Value = 3001
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 3752.5
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 11257.5
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -16157.5

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0

$ws.Range("H132").Value = 3323.6365
$ws.Range("I132").Value = 3166.9285
$ws.Range("K132").Value = 9500.7855
$ws.Range("M132").Value = -6970.7855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 6057
$ws.Range("I8").Value = 7839.8
$ws.Range("J8").Value = 1600
$ws.Range("K8").Value = 7839.8
$ws.Range("L8").Value = 1600
$ws.Range("M8").Value = -7699.8
$ws.Range("N8").Value = -1880

$ws.Range("H86").Value = 27500
$ws.Range("J86").Value = 27500
$ws.Range("L86").Value = 27500
$ws.Range("N86").Value = -29746

$ws.Range("H89").Value = 27500
$ws.Range("J89").Value = 27500
$ws.Range("L89").Value = 137500
$ws.Range("N89").Value = -148732

$ws.Range("H122").Value = 3374.682
$ws.Range("I122").Value = 2028.4375
$ws.Range("K122").Value = 6085.3125
$ws.Range("M122").Value = -3635.3125

$ws.Range("H128").Value = 137185.25
$ws.Range("J128").Value = 137185.25
$ws.Range("L128").Value = 137185.25
$ws.Range("N128").Value = -147145.25

$ws.Range("H132").Value = 3421.2285
$ws.Range("I132").Value = 3027.8147
$ws.Range("J132").Value = 4749
$ws.Range("K132").Value = 9083.444100000001
$ws.Range("L132").Value = 14247
$ws.Range("M132").Value = -6553.444100000001
$ws.Range("N132").Value = -19307
